$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.843.65'
$ws.Range("E2").Value = '  -1.29%  '
$ws.Range("D3").Value = '2.677.92'
$ws.Range("E3").Value = '  -1.21%  '
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '600.08'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.00%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.37'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.67%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.546'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.40%  '
$ws.Range("D9").Value = '2.677.09'
$ws.Range("E9").Value = '  -1.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.145'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +2.91%  '
$ws.Range("E11").Value = '  +1.22%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.364'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +0.82%  '
$ws.Range("E13").Value = '  -1.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.98'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -1.68%  '
$ws.Range("D15").Value = '3.164.26'
$ws.Range("E15").Value = '  -0.93%  '
$ws.Range("E16").Value = '  -2.11%  '
$ws.Range("D17").Value = '67.726.92'
$ws.Range("E17").Value = '  -1.38%  '
$ws.Range("D18").Value = '2.677.32'
$ws.Range("E18").Value = '  -0.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.78'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.79'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.36%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '364.83'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  -3.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.85'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.05'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -3.58%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '71.06'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -4.00%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.19'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.73%  '
$ws.Range("D28").Value = '2.816.63'
$ws.Range("E28").Value = '  -0.83%  '
$ws.Range("E29").Value = '  -2.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '559.64'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.05'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.41'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -3.31%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.94'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.69%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.131'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("E36").Value = '  -0.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.57'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.56'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '156.09'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -2.80%  '
$ws.Range("E40").Value = '  -1.59%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.34'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.84'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -4.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.96'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.54'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -6.06%  '
$ws.Range("E45").Value = '  +0.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.31'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.01%  '
$ws.Range("D47").Value = '0.0₆0302'
$ws.Range("E47").Value = '  -5.00%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.593'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '153.94'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -2.60%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.88'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -1.68%  '
$ws.Range("E51").Value = '  -3.14%  '